$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell C1 - same style as A1/B1 (bold, centered, bordered)
$ws.Range("C1").Value = "2025-03-30"
$ws.Range("C1").Style = $ws.Range("B1").Style

# Row 2: C2 stays blank (inlineStr placeholder in diff - no value entered)
$ws.Range("C2").Value = ""

# Row 3: C3 gets a check-in time
$ws.Range("C3").Value = "✅ 21:58"

# New row 4: Hema, blank B4, check-in time in C4
$ws.Range("A4").Value = "Hema"
$ws.Range("B4").Value = ""
$ws.Range("C4").Value = "✅ 21:58"
